$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "check group" header row for steel members, marked with a red "* New" label
$ws.Range("A26").Value = "* New"
$ws.Range("A26").Font.Color = 255

# New constraint rows for the additional steel-member check groups
$ws.Range("A27").Value = 0
$ws.Range("B27").Value = -6
$ws.Range("C27").Value = -2
$ws.Range("D27").Value = "111111"

$ws.Range("A28").Value = 12
$ws.Range("B28").Value = -6
$ws.Range("C28").Value = -2
$ws.Range("D28").Value = "111111"

$ws.Range("A29").Value = 21
$ws.Range("B29").Value = -6
$ws.Range("C29").Value = -2
$ws.Range("D29").Value = "111111"

$ws.Range("A30").Value = 68
$ws.Range("B30").Value = -4
$ws.Range("C30").Value = -2
$ws.Range("D30").Value = "111111"

$ws.Range("F28").Select()
